# The sheet "Sheet2" (doc_type master data) had the document-type rows for
# DOC011 .. DOC017 (Contrat de location .. Facture d'eau) removed from the
# database, so the exported sheet no longer contains them. DOC018 (Carte
# d'assurance), formerly the last row (19), shifts up to become the new
# last row (12).
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Rows 12-18 held DOC011 through DOC017 - delete them as whole rows so the
# remaining rows (DOC018 etc.) shift up and the used range shrinks from
# A1:E19 to A1:E12.
$ws.Range("A12:E18").EntireRow.Delete()

# Reflect the manual-edit session's final cursor position/selection.
$ws.Activate()
$ws.Range("C22").Select()
